# Commit message: "new update to location"
#
# This adds a new "Notes" worksheet after the existing
# "SF Golden Retrievers Docs" sheet, makes it the active/selected tab,
# puts a single note in cell A1, and leaves the original sheet's
# selection parked on A2 (no longer the selected tab).

$wb = $excel.ActiveWorkbook
$docsSheet = $wb.Worksheets.Item(1)

# Move the selection on the original sheet off of D3 and onto A2
# before we add/activate the new sheet (matches the diff's new
# <selection activeCell="A2" sqref="A2"/> on sheet1, with tabSelected
# no longer present there).
$docsSheet.Range("A2").Select()

# Insert the new "Notes" sheet right after the docs sheet so the sheet
# order becomes: "SF Golden Retrievers Docs", "Notes".
$notesSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $docsSheet)
$notesSheet.Name = "Notes"

# Populate the new sheet and make sure it ends up the active tab.
$notesSheet.Range("A1").Value = "Need to add nav button"
$notesSheet.Activate()
